# Add a new "2022-Q3" quarter sheet and a matching summary row.
#
# 1) Insert a brand-new worksheet named "2022-Q3" right before the existing
#    "2022-Q2" sheet and fill it with the Q3 fund-holding table.
# 2) Insert a new row into the "总计" (summary) sheet for 2022-Q3 and shift
#    the pre-existing quarters down by one row.
#
# Note: once Worksheets.Add() runs, any worksheet object obtained *before*
# that call becomes stale, so sheets are re-fetched by name afterwards.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) New "2022-Q3" worksheet, inserted just before "2022-Q2"
# ---------------------------------------------------------------------------

$q2Before = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($q2Before)
$q3.Name = "2022-Q3"

# Re-fetch "2022-Q2" -- the reference obtained before Add() is stale now.
$q2 = $wb.Worksheets.Item("2022-Q2")

# Header row text.
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Copy the bold/bordered header formatting from the existing quarter sheet.
$q2.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

# Data rows: (index, fund code, fund name, fund size, stock position,
# position ratio, holding value (billion), position rank)
$q3Data = @(
    @(0, "270002", "广发稳健增长混合A",                 "160.46", "38.78", "1.89", "3.0327", 8),
    @(1, "100020", "富国天益价值混合A",                 "60.76",  "89.95", "3.04", "1.8471", 10),
    @(2, "160926", "大成创业板两年定期开放混合A",         "7.68",   "64.14", "3.09", "0.2373", 7),
    @(3, "501070", "广发睿阳三年定期开放混合",            "6.62",   "51.01", "3.54", "0.2343", 6),
    @(4, "011078", "诺德品质消费6个月持有期混合",         "3.86",   "93.34", "3.58", "0.1382", 4),
    @(5, "009798", "大成创业板两年定期开放混合C",         "2.71",   "64.14", "3.09", "0.0837", 7),
    @(6, "001742", "广发百发大数据策略精选灵活配置混合E",  "2.51",   "40.85", "2.61", "0.0655", 7),
    @(7, "009326", "广发稳健增长混合C",                  "2.30",   "38.78", "1.89", "0.0435", 8),
    @(8, "011307", "富国天益价值混合C",                  "0.49",   "89.95", "3.04", "0.0149", 10),
    @(9, "001741", "广发百发大数据策略精选灵活配置混合A",  "0.21",   "40.85", "2.61", "0.0055", 7)
)

$row = 2
foreach ($rec in $q3Data) {
    $q3.Cells.Item($row, 1).Value = $rec[0]

    # Fund code, size, position%, ratio and holding value all look numeric
    # ("270002", "160.46", ...) so a plain .Value assignment would silently
    # turn them into numbers. Prefix with an apostrophe to force text entry
    # (just like typing it into Excel), then reset the style to "Normal" so
    # the quote-prefix formatting flag doesn't linger on the cell.
    $q3.Cells.Item($row, 2).Value = "'" + $rec[1]
    $q3.Cells.Item($row, 2).Style = "Normal"

    $q3.Cells.Item($row, 3).Value = $rec[2]

    $q3.Cells.Item($row, 4).Value = "'" + $rec[3]
    $q3.Cells.Item($row, 4).Style = "Normal"

    $q3.Cells.Item($row, 5).Value = "'" + $rec[4]
    $q3.Cells.Item($row, 5).Style = "Normal"

    $q3.Cells.Item($row, 6).Value = "'" + $rec[5]
    $q3.Cells.Item($row, 6).Style = "Normal"

    $q3.Cells.Item($row, 7).Value = "'" + $rec[6]
    $q3.Cells.Item($row, 7).Style = "Normal"

    $q3.Cells.Item($row, 8).Value = $rec[7]

    $row = $row + 1
}

# Column-A index cells use the bold/bordered style too.
$q2.Range("A2").Copy()
$q3.Range("A2:A11").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) "总计" summary sheet: insert 2022-Q3 row, push older quarters down
# ---------------------------------------------------------------------------

$totalWs = $wb.Worksheets.Item("总计")

# Shift existing rows 2..6 down to 3..7 (bottom-up so nothing is clobbered).
$totalWs.Cells.Item(7, 1).Value = 5
$totalWs.Cells.Item(7, 2).Value = "2021-Q2"
$totalWs.Cells.Item(7, 3).Value = 29
$totalWs.Cells.Item(7, 4).Value = 13.33

$totalWs.Cells.Item(6, 1).Value = 4
$totalWs.Cells.Item(6, 2).Value = "2021-Q3"
$totalWs.Cells.Item(6, 3).Value = 45
$totalWs.Cells.Item(6, 4).Value = 23.75

$totalWs.Cells.Item(5, 1).Value = 3
$totalWs.Cells.Item(5, 2).Value = "2021-Q4"
$totalWs.Cells.Item(5, 3).Value = 43
$totalWs.Cells.Item(5, 4).Value = 31.89

$totalWs.Cells.Item(4, 1).Value = 2
$totalWs.Cells.Item(4, 2).Value = "2022-Q1"
$totalWs.Cells.Item(4, 3).Value = 31
$totalWs.Cells.Item(4, 4).Value = 22.65

$totalWs.Cells.Item(3, 1).Value = 1
$totalWs.Cells.Item(3, 2).Value = "2022-Q2"
$totalWs.Cells.Item(3, 3).Value = 29
$totalWs.Cells.Item(3, 4).Value = 21.85

$totalWs.Cells.Item(2, 1).Value = 0
$totalWs.Cells.Item(2, 2).Value = "2022-Q3"
$totalWs.Cells.Item(2, 3).Value = 10
$totalWs.Cells.Item(2, 4).Value = 5.7

# Re-apply the bold/bordered "index column" style to the newly-created A7
# cell (it has no pre-existing formatting to inherit).
$totalWs.Range("A6").Copy()
$totalWs.Range("A7").PasteSpecial(-4122)
